$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.051203666666667
$ws.Range("H2").Value = 18.153611
$ws.Range("I2").Value = 0.07608037240065801
$ws.Range("J2").Value = 0.07775008964215516
$ws.Range("M2").Value = 9.363528666666667
$ws.Range("N2").Value = 28.090586
$ws.Range("O2").Value = 0.04175743631338733
$ws.Range("P2").Value = 0.04324026421082073
$ws.Range("Q2").Value = 56.66061900067179
$ws.Range("R2").Value = 509.9455710060461
$ws.Range("S2").Value = 0.003176921305219268
$ws.Range("T2").Value = 0.003361934418541785

# Row 3
$ws.Range("G3").Value = 6.051203666666667
$ws.Range("H3").Value = 18.153611
$ws.Range("I3").Value = 0.07608037240065801
$ws.Range("J3").Value = 0.07775008964215516
$ws.Range("O3").Value = 0.1749266505387075
$ws.Range("P3").Value = 0.1811383852696593
$ws.Range("Q3").Value = 237.3577780219178
$ws.Range("R3").Value = 2136.22000219726
$ws.Range("S3").Value = 0.01330848471578463
$ws.Range("T3").Value = 0.01408352569235125

# Row 4
$ws.Range("G4").Value = 6.051203666666667
$ws.Range("H4").Value = 18.153611
$ws.Range("I4").Value = 0.07608037240065801
$ws.Range("J4").Value = 0.07775008964215516
$ws.Range("M4").Value = 71.284935
$ws.Range("N4").Value = 213.854805
$ws.Range("O4").Value = 0.3179011075133629
$ws.Range("P4").Value = 0.3291899382573772
$ws.Range("Q4").Value = 431.359660050095
$ws.Range("R4").Value = 3882.236940450855
$ws.Range("S4").Value = 0.02418603464619827
$ws.Range("T4").Value = 0.0255945472088066

# Row 5
$ws.Range("G5").Value = 6.051203666666667
$ws.Range("H5").Value = 18.153611
$ws.Range("I5").Value = 0.07608037240065801
$ws.Range("J5").Value = 0.07775008964215516
$ws.Range("M5").Value = 23.0690325
$ws.Range("N5").Value = 46.138065
$ws.Range("O5").Value = 0.1028782726814826
$ws.Range("P5").Value = 0.07102102180339065
$ws.Range("Q5").Value = 139.5954140504525
$ws.Range("R5").Value = 837.572484302715
$ws.Range("S5").Value = 0.007827017297543638
$ws.Range("T5").Value = 0.005521890811691079

# Row 6
$ws.Range("G6").Value = 6.051203666666667
$ws.Range("H6").Value = 18.153611
$ws.Range("I6").Value = 0.07608037240065801
$ws.Range("J6").Value = 0.07775008964215516
$ws.Range("M6").Value = 81.293813
$ws.Range("N6").Value = 243.881439
$ws.Range("O6").Value = 0.3625365329530597
$ws.Range("P6").Value = 0.3754103904587522
$ws.Range("Q6").Value = 491.9254193029143
$ws.Range("R6").Value = 4427.328773726229
$ws.Range("S6").Value = 0.02758191443591221
$ws.Range("T6").Value = 0.02918819151076445

# Row 7
$ws.Range("I7").Value = 0.2215826302097334
$ws.Range("J7").Value = 0.2264456497560767
$ws.Range("M7").Value = 9.363528666666667
$ws.Range("N7").Value = 28.090586
$ws.Range("O7").Value = 0.04175743631338733
$ws.Range("P7").Value = 0.04324026421082073
$ws.Range("Q7").Value = 165.0229696742634
$ws.Range("R7").Value = 1485.20672706837
$ws.Range("S7").Value = 0.009252722569135797
$ws.Range("T7").Value = 0.009791569724843728

# Row 8
$ws.Range("I8").Value = 0.2215826302097334
$ws.Range("J8").Value = 0.2264456497560767
$ws.Range("O8").Value = 0.1749266505387075
$ws.Range("P8").Value = 0.1811383852696593
$ws.Range("S8").Value = 0.03876070732014567
$ws.Range("T8").Value = 0.04101799934815455

# Row 9
$ws.Range("I9").Value = 0.2215826302097334
$ws.Range("J9").Value = 0.2264456497560767
$ws.Range("M9").Value = 71.284935
$ws.Range("N9").Value = 213.854805
$ws.Range("O9").Value = 0.3179011075133629
$ws.Range("P9").Value = 0.3291899382573772
$ws.Range("Q9").Value = 1256.326763714025
$ws.Range("R9").Value = 11306.94087342623
$ws.Range("S9").Value = 0.07044136354939819
$ws.Range("T9").Value = 0.07454362946185455

# Row 10
$ws.Range("I10").Value = 0.2215826302097334
$ws.Range("J10").Value = 0.2264456497560767
$ws.Range("M10").Value = 23.0690325
$ws.Range("N10").Value = 46.138065
$ws.Range("O10").Value = 0.1028782726814826
$ws.Range("P10").Value = 0.07102102180339065
$ws.Range("Q10").Value = 406.5689748154875
$ws.Range("R10").Value = 2439.413848892925
$ws.Range("S10").Value = 0.02279603825219707
$ws.Range("T10").Value = 0.01608240142860929

# Row 11
$ws.Range("I11").Value = 0.2215826302097334
$ws.Range("J11").Value = 0.2264456497560767
$ws.Range("M11").Value = 81.293813
$ws.Range("N11").Value = 243.881439
$ws.Range("O11").Value = 0.3625365329530597
$ws.Range("P11").Value = 0.3754103904587522
$ws.Range("Q11").Value = 1432.723379719195
$ws.Range("R11").Value = 12894.51041747275
$ws.Range("S11").Value = 0.08033179851885665
$ws.Range("T11").Value = 0.08501004979261459

# Row 12
$ws.Range("G12").Value = 33.62840566666667
$ws.Range("H12").Value = 100.885217
$ws.Range("I12").Value = 0.4228021014155913
$ws.Range("J12").Value = 0.432081235260482
$ws.Range("M12").Value = 9.363528666666667
$ws.Range("N12").Value = 28.090586
$ws.Range("O12").Value = 0.04175743631338733
$ws.Range("P12").Value = 0.04324026421082073
$ws.Range("Q12").Value = 314.8805404741291
$ws.Range("R12").Value = 2833.924864267162
$ws.Range("S12").Value = 0.01765513182302789
$ws.Range("T12").Value = 0.01868330677320103

# Row 13
$ws.Range("G13").Value = 33.62840566666667
$ws.Range("H13").Value = 100.885217
$ws.Range("I13").Value = 0.4228021014155913
$ws.Range("J13").Value = 0.432081235260482
$ws.Range("O13").Value = 0.1749266505387075
$ws.Range("P13").Value = 0.1811383852696593
$ws.Range("Q13").Value = 1319.070401055691
$ws.Range("R13").Value = 11871.63360950122
$ws.Range("S13").Value = 0.07395935544135629
$ws.Range("T13").Value = 0.07826649726040348

# Row 14
$ws.Range("G14").Value = 33.62840566666667
$ws.Range("H14").Value = 100.885217
$ws.Range("I14").Value = 0.4228021014155913
$ws.Range("J14").Value = 0.432081235260482
$ws.Range("M14").Value = 71.284935
$ws.Range("N14").Value = 213.854805
$ws.Range("O14").Value = 0.3179011075133629
$ws.Range("P14").Value = 0.3291899382573772
$ws.Range("Q14").Value = 2397.198712101965
$ws.Range("R14").Value = 21574.78840891768
$ws.Range("S14").Value = 0.1344092562989937
$ws.Range("T14").Value = 0.1422367951575693

# Row 15
$ws.Range("G15").Value = 33.62840566666667
$ws.Range("H15").Value = 100.885217
$ws.Range("I15").Value = 0.4228021014155913
$ws.Range("J15").Value = 0.432081235260482
$ws.Range("M15").Value = 23.0690325
$ws.Range("N15").Value = 46.138065
$ws.Range("O15").Value = 0.1028782726814826
$ws.Range("P15").Value = 0.07102102180339065
$ws.Range("Q15").Value = 775.7747832475175
$ws.Range("R15").Value = 4654.648699485105
$ws.Range("S15").Value = 0.04349714987973707
$ws.Range("T15").Value = 0.03068685083027066

# Row 16
$ws.Range("G16").Value = 33.62840566666667
$ws.Range("H16").Value = 100.885217
$ws.Range("I16").Value = 0.4228021014155913
$ws.Range("J16").Value = 0.432081235260482
$ws.Range("M16").Value = 81.293813
$ws.Range("N16").Value = 243.881439
$ws.Range("O16").Value = 0.3625365329530597
$ws.Range("P16").Value = 0.3754103904587522
$ws.Range("Q16").Value = 2733.78132175414
$ws.Range("R16").Value = 24604.03189578726
$ws.Range("S16").Value = 0.1532812079724764
$ws.Range("T16").Value = 0.1622077852390375

# Row 17
$ws.Range("G17").Value = 5.124275
$ws.Range("H17").Value = 10.24855
$ws.Range("I17").Value = 0.06442631445887793
$ws.Range("J17").Value = 0.04389350863594627
$ws.Range("M17").Value = 9.363528666666667
$ws.Range("N17").Value = 28.090586
$ws.Range("O17").Value = 0.04175743631338733
$ws.Range("P17").Value = 0.04324026421082073
$ws.Range("Q17").Value = 47.98129585838333
$ws.Range("R17").Value = 287.8877751503
$ws.Range("S17").Value = 0.002690277722922861
$ws.Range("T17").Value = 0.001897966910558258

# Row 18
$ws.Range("G18").Value = 5.124275
$ws.Range("H18").Value = 10.24855
$ws.Range("I18").Value = 0.06442631445887793
$ws.Range("J18").Value = 0.04389350863594627
$ws.Range("O18").Value = 0.1749266505387075
$ws.Range("P18").Value = 0.1811383852696593
$ws.Range("Q18").Value = 200.9991061238333
$ws.Range("R18").Value = 1205.994636743
$ws.Range("S18").Value = 0.01126987939484502
$ws.Range("T18").Value = 0.007950799278135152

# Row 19
$ws.Range("G19").Value = 5.124275
$ws.Range("H19").Value = 10.24855
$ws.Range("I19").Value = 0.06442631445887793
$ws.Range("J19").Value = 0.04389350863594627
$ws.Range("M19").Value = 71.284935
$ws.Range("N19").Value = 213.854805
$ws.Range("O19").Value = 0.3179011075133629
$ws.Range("P19").Value = 0.3291899382573772
$ws.Range("Q19").Value = 365.283610297125
$ws.Range("R19").Value = 2191.70166178275
$ws.Range("S19").Value = 0.02048119671948148
$ws.Range("T19").Value = 0.0144493013977668

# Row 20
$ws.Range("G20").Value = 5.124275
$ws.Range("H20").Value = 10.24855
$ws.Range("I20").Value = 0.06442631445887793
$ws.Range("J20").Value = 0.04389350863594627
$ws.Range("M20").Value = 23.0690325
$ws.Range("N20").Value = 46.138065
$ws.Range("O20").Value = 0.1028782726814826
$ws.Range("P20").Value = 0.07102102180339065
$ws.Range("Q20").Value = 118.2120665139375
$ws.Range("R20").Value = 472.8482660557499
$ws.Range("S20").Value = 0.00662806794676339
$ws.Range("T20").Value = 0.003117361833860855

# Row 21
$ws.Range("G21").Value = 5.124275
$ws.Range("H21").Value = 10.24855
$ws.Range("I21").Value = 0.06442631445887793
$ws.Range("J21").Value = 0.04389350863594627
$ws.Range("M21").Value = 81.293813
$ws.Range("N21").Value = 243.881439
$ws.Range("O21").Value = 0.3625365329530597
$ws.Range("P21").Value = 0.3754103904587522
$ws.Range("Q21").Value = 416.571853610575
$ws.Range("R21").Value = 2499.43112166345
$ws.Range("S21").Value = 0.02335689267486519
$ws.Range("T21").Value = 0.0164780792156252

# Row 22
$ws.Range("G22").Value = 17.10908866666667
$ws.Range("H22").Value = 51.327266
$ws.Range("I22").Value = 0.2151085815151395
$ws.Range("J22").Value = 0.2198295167053399
$ws.Range("M22").Value = 9.363528666666667
$ws.Range("N22").Value = 28.090586
$ws.Range("O22").Value = 0.04175743631338733
$ws.Range("P22").Value = 0.04324026421082073
$ws.Range("Q22").Value = 160.2014421908751
$ws.Range("R22").Value = 1441.812979717876
$ws.Range("S22").Value = 0.008982382893081524
$ws.Range("T22").Value = 0.009505486383675927

# Row 23
$ws.Range("G23").Value = 17.10908866666667
$ws.Range("H23").Value = 51.327266
$ws.Range("I23").Value = 0.2151085815151395
$ws.Range("J23").Value = 0.2198295167053399
$ws.Range("O23").Value = 0.1749266505387075
$ws.Range("P23").Value = 0.1811383852696593
$ws.Range("Q23").Value = 671.1020639199512
$ws.Range("R23").Value = 6039.91857527956
$ws.Range("S23").Value = 0.03762822366657587
$ws.Range("T23").Value = 0.03981956369061486

# Row 24
$ws.Range("G24").Value = 17.10908866666667
$ws.Range("H24").Value = 51.327266
$ws.Range("I24").Value = 0.2151085815151395
$ws.Range("J24").Value = 0.2198295167053399
$ws.Range("M24").Value = 71.284935
$ws.Range("N24").Value = 213.854805
$ws.Range("O24").Value = 0.3179011075133629
$ws.Range("P24").Value = 0.3291899382573772
$ws.Range("Q24").Value = 1219.62027351257
$ws.Range("R24").Value = 10976.58246161313
$ws.Range("S24").Value = 0.06838325629929135
$ws.Range("T24").Value = 0.07236566503137991

# Row 25
$ws.Range("G25").Value = 17.10908866666667
$ws.Range("H25").Value = 51.327266
$ws.Range("I25").Value = 0.2151085815151395
$ws.Range("J25").Value = 0.2198295167053399
$ws.Range("M25").Value = 23.0690325
$ws.Range("N25").Value = 46.138065
$ws.Range("O25").Value = 0.1028782726814826
$ws.Range("P25").Value = 0.07102102180339065
$ws.Range("Q25").Value = 394.690122496715
$ws.Range("R25").Value = 2368.14073498029
$ws.Range("S25").Value = 0.02212999930524145
$ws.Range("T25").Value = 0.01561251689895877

# Row 26
$ws.Range("G26").Value = 17.10908866666667
$ws.Range("H26").Value = 51.327266
$ws.Range("I26").Value = 0.2151085815151395
$ws.Range("J26").Value = 0.2198295167053399
$ws.Range("M26").Value = 81.293813
$ws.Range("N26").Value = 243.881439
$ws.Range("O26").Value = 0.3625365329530597
$ws.Range("P26").Value = 0.3754103904587522
$ws.Range("Q26").Value = 1390.863054668419
$ws.Range("R26").Value = 12517.76749201578
$ws.Range("S26").Value = 0.07798471935094929
$ws.Range("T26").Value = 0.08252628470071044
